$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting (style index reused) for new header columns AK:AR from the existing FGSM block (AC:AJ) ---
$ws.Range("AC1:AJ1").Copy() | Out-Null
$ws.Range("AK1:AR1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("AC2:AJ2").Copy() | Out-Null
$ws.Range("AK2:AR2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 1 header: BOUNDARY attack block label, then merge AK1:AR1 ---
$ws.Range("AK1").Value = "BOUNDARY"
$ws.Range("AK1:AR1").Merge() | Out-Null

# --- Row 2 epsilon labels (stored as text, matching the existing ε row) ---
$ws.Range("AK2").Value = "0.01"
$ws.Range("AL2").Value = "0.02"
$ws.Range("AM2").Value = "0.03"
$ws.Range("AN2").Value = "0.04"
$ws.Range("AO2").Value = "0.05"
$ws.Range("AP2").Value = "0.07"
$ws.Range("AQ2").Value = "0.10"
$ws.Range("AR2").Value = "0.20"

# --- New BOUNDARY-attack data columns (AK:AR) for every metric row ---
# Row 4
$ws.Range("AK4").Value = 383.1023151906331
$ws.Range("AL4").Value = 386.3325114313761
$ws.Range("AM4").Value = 392.373773021698
$ws.Range("AN4").Value = 391.0664524650574
$ws.Range("AO4").Value = 401.8060641860962
$ws.Range("AP4").Value = 405.8325726763408
$ws.Range("AQ4").Value = 465.6313930575053
$ws.Range("AR4").Value = 564.4294561258952

# Row 5
$ws.Range("AK5").Value = 490.2629962789269
$ws.Range("AL5").Value = 495.4276808743702
$ws.Range("AM5").Value = 502.0055808490563
$ws.Range("AN5").Value = 499.1427309591197
$ws.Range("AO5").Value = 521.6158783918642
$ws.Range("AP5").Value = 509.0743864712236
$ws.Range("AQ5").Value = 592.0152907778793
$ws.Range("AR5").Value = 721.2074289555613

# Row 6
$ws.Range("AK6").Value = 0.9991674931885181
$ws.Range("AL6").Value = 0.9991509824837742
$ws.Range("AM6").Value = 0.9991365178272575
$ws.Range("AN6").Value = 0.9991244010732634
$ws.Range("AO6").Value = 0.9990247489845747
$ws.Range("AP6").Value = 0.9990994162118598
$ws.Range("AQ6").Value = 0.9987554554869923
$ws.Range("AR6").Value = 0.998162011269959

# Row 7
$ws.Range("AK7").Value = 370.3934350649516
$ws.Range("AL7").Value = 372.6032810401916
$ws.Range("AM7").Value = 379.8744000562032
$ws.Range("AN7").Value = 387.8404699516296
$ws.Range("AO7").Value = 390.8763733609517
$ws.Range("AP7").Value = 458.6322962760925
$ws.Range("AQ7").Value = 472.8261692492167
$ws.Range("AR7").Value = 665.7321542930603

# Row 8
$ws.Range("AK8").Value = 484.8727072655358
$ws.Range("AL8").Value = 489.7849810900284
$ws.Range("AM8").Value = 496.964290138574
$ws.Range("AN8").Value = 497.0586226283393
$ws.Range("AO8").Value = 500.7750292329047
$ws.Range("AP8").Value = 574.6016949385179
$ws.Range("AQ8").Value = 600.4894722591896
$ws.Range("AR8").Value = 831.3046039408679

# Row 9
$ws.Range("AK9").Value = 0.9991923532698864
$ws.Range("AL9").Value = 0.9991651103922494
$ws.Range("AM9").Value = 0.9991402639451001
$ws.Range("AN9").Value = 0.9991474490417931
$ws.Range("AO9").Value = 0.9991561364028338
$ws.Range("AP9").Value = 0.9988542530996758
$ws.Range("AQ9").Value = 0.9987415250533545
$ws.Range("AR9").Value = 0.9973933839876078

# Row 10
$ws.Range("AK10").Value = 305.3789920997619
$ws.Range("AL10").Value = 312.5076346143087
$ws.Range("AM10").Value = 318.2210095659892
$ws.Range("AN10").Value = 329.9936900393168
$ws.Range("AO10").Value = 355.7949362881978
$ws.Range("AP10").Value = 408.0614822705587
$ws.Range("AQ10").Value = 413.5686668078105
$ws.Range("AR10").Value = 643.0573875713349

# Row 11
$ws.Range("AK11").Value = 416.8893104573613
$ws.Range("AL11").Value = 426.0453997390265
$ws.Range("AM11").Value = 429.2889126112187
$ws.Range("AN11").Value = 427.7417020012782
$ws.Range("AO11").Value = 478.4552467500266
$ws.Range("AP11").Value = 525.4992721913908
$ws.Range("AQ11").Value = 563.9345859192189
$ws.Range("AR11").Value = 817.6874174168767

# Row 12
$ws.Range("AK12").Value = 0.9993484402843759
$ws.Range("AL12").Value = 0.9993132180291499
$ws.Range("AM12").Value = 0.9993078200444149
$ws.Range("AN12").Value = 0.9993208202900299
$ws.Range("AO12").Value = 0.9991477585313444
$ws.Range("AP12").Value = 0.9989757896469301
$ws.Range("AQ12").Value = 0.9987403941159547
$ws.Range("AR12").Value = 0.9973520020895164

# --- Tiny floating-point refresh on a handful of pre-existing SIM cells (rows 6, 9, 12) ---
$ws.Range("AC6").Value = 0.9988626594988415
$ws.Range("AG6").Value = 0.9967461276332034
$ws.Range("AH6").Value = 0.9950072733341898
$ws.Range("AI6").Value = 0.9924597667458618
$ws.Range("H6").Value = 0.9990505464682528
$ws.Range("K6").Value = 0.9985619634325859
$ws.Range("T6").Value = 0.9975501151841627
$ws.Range("AC9").Value = 0.9989831240878908
$ws.Range("AE9").Value = 0.9983964271126717
$ws.Range("D9").Value = 0.9991508556045917
$ws.Range("E9").Value = 0.9991506278049966
$ws.Range("L9").Value = 0.9952898719064709
$ws.Range("M9").Value = 0.9991506278049966
$ws.Range("O9").Value = 0.9990136185314769
$ws.Range("Q9").Value = 0.9988283505152307
$ws.Range("R9").Value = 0.9985745158613824
$ws.Range("U9").Value = 0.9990191256829541
$ws.Range("Z9").Value = 0.9967828838849236
$ws.Range("AB12").Value = 0.9868612077339067
$ws.Range("AD12").Value = 0.9987828975885487
$ws.Range("AG12").Value = 0.9972041427377818
$ws.Range("AI12").Value = 0.9923753191835871
$ws.Range("C12").Value = 0.9993617696409695
$ws.Range("F12").Value = 0.9992456130615981
$ws.Range("H12").Value = 0.9992331020769409
$ws.Range("K12").Value = 0.9988721190136414
$ws.Range("M12").Value = 0.999383747802967
$ws.Range("N12").Value = 0.9993956324437535
$ws.Range("O12").Value = 0.9993946161303718
$ws.Range("P12").Value = 0.9993770742595574
$ws.Range("R12").Value = 0.9992429469063333
$ws.Range("T12").Value = 0.9974327371754078
$ws.Range("Z12").Value = 0.9967571849739728

Write-Host "BOUNDARY attack columns (AK:AR) added for SEED 314"
